# Natmi LR-pair output update (Clcf1-Lifr, YoungD7) following Dr Hou advice.
# Sending-cluster set now includes "ECs" alongside FAPs/M2/sCs, expanding the
# 4x3 Sending x Target matrix (rows 2-13) into a 4x4 matrix (rows 2-17).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = New-Object 'object[,]' 16,20

# Row 2: ECs -> ECs
$values[0,0] = "ECs"
$values[0,1] = "Clcf1"
$values[0,2] = "Lifr"
$values[0,3] = "ECs"
$values[0,4] = 2
$values[0,5] = 0.6666666666666666
$values[0,6] = 1.679012
$values[0,7] = 5.037036000000001
$values[0,8] = 0.1178149724053671
$values[0,9] = 0.1178149724053671
$values[0,10] = 3
$values[0,11] = 1
$values[0,12] = 26.31197333333334
$values[0,13] = 78.93592000000001
$values[0,14] = 0.2261559208386891
$values[0,15] = 0.2261559208386891
$values[0,16] = 44.17811897034668
$values[0,17] = 397.6030707331201
$values[0,18] = 0.02664455357292056
$values[0,19] = 0.02664455357292055

# Row 3: ECs -> FAPs
$values[1,0] = "ECs"
$values[1,1] = "Clcf1"
$values[1,2] = "Lifr"
$values[1,3] = "FAPs"
$values[1,4] = 2
$values[1,5] = 0.6666666666666666
$values[1,6] = 1.679012
$values[1,7] = 5.037036000000001
$values[1,8] = 0.1178149724053671
$values[1,9] = 0.1178149724053671
$values[1,10] = 3
$values[1,11] = 1
$values[1,12] = 53.74150833333334
$values[1,13] = 161.224525
$values[1,14] = 0.461917475759518
$values[1,15] = 0.461917475759518
$values[1,16] = 90.23263738976668
$values[1,17] = 812.0937365079001
$values[1,18] = 0.05442079466016447
$values[1,19] = 0.05442079466016446

# Row 4: ECs -> M2
$values[2,0] = "ECs"
$values[2,1] = "Clcf1"
$values[2,2] = "Lifr"
$values[2,3] = "M2"
$values[2,4] = 2
$values[2,5] = 0.6666666666666666
$values[2,6] = 1.679012
$values[2,7] = 5.037036000000001
$values[2,8] = 0.1178149724053671
$values[2,9] = 0.1178149724053671
$values[2,10] = 3
$values[2,11] = 1
$values[2,12] = 21.978693
$values[2,13] = 65.93607899999999
$values[2,14] = 0.1889106336220259
$values[2,15] = 0.1889106336220259
$values[2,16] = 36.90248929131599
$values[2,17] = 332.122403621844
$values[2,18] = 0.02225650108725941
$values[2,19] = 0.02225650108725941

# Row 5: ECs -> sCs
$values[3,0] = "ECs"
$values[3,1] = "Clcf1"
$values[3,2] = "Lifr"
$values[3,3] = "sCs"
$values[3,4] = 2
$values[3,5] = 0.6666666666666666
$values[3,6] = 1.679012
$values[3,7] = 5.037036000000001
$values[3,8] = 0.1178149724053671
$values[3,9] = 0.1178149724053671
$values[3,10] = 3
$values[3,11] = 1
$values[3,12] = 14.31221833333333
$values[3,13] = 42.936655
$values[3,14] = 0.123015969779767
$values[3,15] = 0.123015969779767
$values[3,16] = 24.03038632828667
$values[3,17] = 216.27347695458
$values[3,18] = 0.01449312308502273
$values[3,19] = 0.01449312308502273

# Row 6: FAPs -> ECs
$values[4,0] = "FAPs"
$values[4,1] = "Clcf1"
$values[4,2] = "Lifr"
$values[4,3] = "ECs"
$values[4,4] = 3
$values[4,5] = 1
$values[4,6] = 3.197979
$values[4,7] = 9.593937
$values[4,8] = 0.2243997110431275
$values[4,9] = 0.2243997110431275
$values[4,10] = 3
$values[4,11] = 1
$values[4,12] = 26.31197333333334
$values[4,13] = 78.93592000000001
$values[4,14] = 0.2261559208386891
$values[4,15] = 0.2261559208386891
$values[4,16] = 84.14513816856002
$values[4,17] = 757.3062435170401
$values[4,18] = 0.05074932328689426
$values[4,19] = 0.05074932328689425

# Row 7: FAPs -> FAPs
$values[5,0] = "FAPs"
$values[5,1] = "Clcf1"
$values[5,2] = "Lifr"
$values[5,3] = "FAPs"
$values[5,4] = 3
$values[5,5] = 1
$values[5,6] = 3.197979
$values[5,7] = 9.593937
$values[5,8] = 0.2243997110431275
$values[5,9] = 0.2243997110431275
$values[5,10] = 3
$values[5,11] = 1
$values[5,12] = 53.74150833333334
$values[5,13] = 161.224525
$values[5,14] = 0.461917475759518
$values[5,15] = 0.461917475759518
$values[5,16] = 171.864215078325
$values[5,17] = 1546.777935704925
$values[5,18] = 0.1036541480862067
$values[5,19] = 0.1036541480862067

# Row 8: FAPs -> M2
$values[6,0] = "FAPs"
$values[6,1] = "Clcf1"
$values[6,2] = "Lifr"
$values[6,3] = "M2"
$values[6,4] = 3
$values[6,5] = 1
$values[6,6] = 3.197979
$values[6,7] = 9.593937
$values[6,8] = 0.2243997110431275
$values[6,9] = 0.2243997110431275
$values[6,10] = 3
$values[6,11] = 1
$values[6,12] = 21.978693
$values[6,13] = 65.93607899999999
$values[6,14] = 0.1889106336220259
$values[6,15] = 0.1889106336220259
$values[6,16] = 70.28739866144699
$values[6,17] = 632.5865879530229
$values[6,18] = 0.04239149159775674
$values[6,19] = 0.04239149159775674

# Row 9: FAPs -> sCs
$values[7,0] = "FAPs"
$values[7,1] = "Clcf1"
$values[7,2] = "Lifr"
$values[7,3] = "sCs"
$values[7,4] = 3
$values[7,5] = 1
$values[7,6] = 3.197979
$values[7,7] = 9.593937
$values[7,8] = 0.2243997110431275
$values[7,9] = 0.2243997110431275
$values[7,10] = 3
$values[7,11] = 1
$values[7,12] = 14.31221833333333
$values[7,13] = 42.936655
$values[7,14] = 0.123015969779767
$values[7,15] = 0.123015969779767
$values[7,16] = 45.770173673415
$values[7,17] = 411.931563060735
$values[7,18] = 0.02760474807226982
$values[7,19] = 0.02760474807226981

# Row 10: M2 -> ECs
$values[8,0] = "M2"
$values[8,1] = "Clcf1"
$values[8,2] = "Lifr"
$values[8,3] = "ECs"
$values[8,4] = 3
$values[8,5] = 1
$values[8,6] = 1.919382666666667
$values[8,7] = 5.758148
$values[8,8] = 0.134681596027112
$values[8,9] = 0.134681596027112
$values[8,10] = 3
$values[8,11] = 1
$values[8,12] = 26.31197333333334
$values[8,13] = 78.93592000000001
$values[8,14] = 0.2261559208386891
$values[8,15] = 0.2261559208386891
$values[8,16] = 50.50274554179556
$values[8,17] = 454.5247098761601
$values[8,18] = 0.03045904036953584
$values[8,19] = 0.03045904036953584

# Row 11: M2 -> FAPs
$values[9,0] = "M2"
$values[9,1] = "Clcf1"
$values[9,2] = "Lifr"
$values[9,3] = "FAPs"
$values[9,4] = 3
$values[9,5] = 1
$values[9,6] = 1.919382666666667
$values[9,7] = 5.758148
$values[9,8] = 0.134681596027112
$values[9,9] = 0.134681596027112
$values[9,10] = 3
$values[9,11] = 1
$values[9,12] = 53.74150833333334
$values[9,13] = 161.224525
$values[9,14] = 0.461917475759518
$values[9,15] = 0.461917475759518
$values[9,16] = 103.1505195755222
$values[9,17] = 928.3546761797
$values[9,18] = 0.06221178286810669
$values[9,19] = 0.06221178286810668

# Row 12: M2 -> M2
$values[10,0] = "M2"
$values[10,1] = "Clcf1"
$values[10,2] = "Lifr"
$values[10,3] = "M2"
$values[10,4] = 3
$values[10,5] = 1
$values[10,6] = 1.919382666666667
$values[10,7] = 5.758148
$values[10,8] = 0.134681596027112
$values[10,9] = 0.134681596027112
$values[10,10] = 3
$values[10,11] = 1
$values[10,12] = 21.978693
$values[10,13] = 65.93607899999999
$values[10,14] = 0.1889106336220259
$values[10,15] = 0.1889106336220259
$values[10,16] = 42.18552238018799
$values[10,17] = 379.669701421692
$values[10,18] = 0.02544278564270745
$values[10,19] = 0.02544278564270745

# Row 13: M2 -> sCs
$values[11,0] = "M2"
$values[11,1] = "Clcf1"
$values[11,2] = "Lifr"
$values[11,3] = "sCs"
$values[11,4] = 3
$values[11,5] = 1
$values[11,6] = 1.919382666666667
$values[11,7] = 5.758148
$values[11,8] = 0.134681596027112
$values[11,9] = 0.134681596027112
$values[11,10] = 3
$values[11,11] = 1
$values[11,12] = 14.31221833333333
$values[11,13] = 42.936655
$values[11,14] = 0.123015969779767
$values[11,15] = 0.123015969779767
$values[11,16] = 27.47062379054889
$values[11,17] = 247.23561411494
$values[11,18] = 0.01656798714676199
$values[11,19] = 0.01656798714676199

# Row 14: sCs -> ECs
$values[12,0] = "sCs"
$values[12,1] = "Clcf1"
$values[12,2] = "Lifr"
$values[12,3] = "ECs"
$values[12,4] = 3
$values[12,5] = 1
$values[12,6] = 7.454888
$values[12,7] = 22.364664
$values[12,8] = 0.5231037205243934
$values[12,9] = 0.5231037205243932
$values[12,10] = 3
$values[12,11] = 1
$values[12,12] = 26.31197333333334
$values[12,13] = 78.93592000000001
$values[12,14] = 0.2261559208386891
$values[12,15] = 0.2261559208386891
$values[12,16] = 196.1528142589867
$values[12,17] = 1765.37532833088
$values[12,18] = 0.1183030036093385
$values[12,19] = 0.1183030036093384

# Row 15: sCs -> FAPs
$values[13,0] = "sCs"
$values[13,1] = "Clcf1"
$values[13,2] = "Lifr"
$values[13,3] = "FAPs"
$values[13,4] = 3
$values[13,5] = 1
$values[13,6] = 7.454888
$values[13,7] = 22.364664
$values[13,8] = 0.5231037205243934
$values[13,9] = 0.5231037205243932
$values[13,10] = 3
$values[13,11] = 1
$values[13,12] = 53.74150833333334
$values[13,13] = 161.224525
$values[13,14] = 0.461917475759518
$values[13,15] = 0.461917475759518
$values[13,16] = 400.6369255760667
$values[13,17] = 3605.7323301846
$values[13,18] = 0.2416307501450402
$values[13,19] = 0.2416307501450401

# Row 16: sCs -> M2
$values[14,0] = "sCs"
$values[14,1] = "Clcf1"
$values[14,2] = "Lifr"
$values[14,3] = "M2"
$values[14,4] = 3
$values[14,5] = 1
$values[14,6] = 7.454888
$values[14,7] = 22.364664
$values[14,8] = 0.5231037205243934
$values[14,9] = 0.5231037205243932
$values[14,10] = 3
$values[14,11] = 1
$values[14,12] = 21.978693
$values[14,13] = 65.93607899999999
$values[14,14] = 0.1889106336220259
$values[14,15] = 0.1889106336220259
$values[14,16] = 163.848694701384
$values[14,17] = 1474.638252312456
$values[14,18] = 0.09881985529430232
$values[14,19] = 0.09881985529430229

# Row 17: sCs -> sCs
$values[15,0] = "sCs"
$values[15,1] = "Clcf1"
$values[15,2] = "Lifr"
$values[15,3] = "sCs"
$values[15,4] = 3
$values[15,5] = 1
$values[15,6] = 7.454888
$values[15,7] = 22.364664
$values[15,8] = 0.5231037205243934
$values[15,9] = 0.5231037205243932
$values[15,10] = 3
$values[15,11] = 1
$values[15,12] = 14.31221833333333
$values[15,13] = 42.936655
$values[15,14] = 0.123015969779767
$values[15,15] = 0.123015969779767
$values[15,16] = 106.6959847065467
$values[15,17] = 960.2638623589201
$values[15,18] = 0.06435011147571246
$values[15,19] = 0.06435011147571243

# Write the full replacement block in one shot (matches the new A1:T17 dimension).
$ws.Range("A2:T17").Value = $values
